$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "Daftar Nomor Handphone" merged banner ---
$ws.Range("A1:D1").UnMerge()

# --- Drop the old table header row + its two data rows + the trailing blank rows ---
$ws.Rows("3:9").Delete()

# --- Reset formatting across the remaining used area back to the workbook default ---
$ws.Range("A1:D2").Style = "Normal"

# --- Wipe out columns B and D entirely (values + formats); only column A survives ---
$ws.Range("B1:D2").Clear()

# --- New minimal "SMS gateway" phone list: a header and a single number, both as text ---
$ws.Range("A1:A2").NumberFormat = "@"
$ws.Range("A1").Value = "No_HP"
$ws.Range("A2").Value = "081386745521"

# --- Match the column A width from the edited workbook (renders as width 12 in xlsx units) ---
$ws.Columns("A").ColumnWidth = 11.2

# --- Leave the selection where the author's session ended up ---
$null = $ws.Range("D6").Select()
